# Append a new arrivals row (row 8) to the "Main Data" sheet, mirroring the
# existing rows (Number, Date, Time, Flight, From, Short, Airline, Model,
# Aircraft ID, Status, blank, Difference, blank).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Sunday, Jan 08"
$ws.Range("C8").Value = "9:10 PM"
$ws.Range("D8").Value = "FR6639"
$ws.Range("E8").Value = "London"
$ws.Range("F8").Value = "(LTN)"
$ws.Range("G8").Value = "Ryanair "
$ws.Range("H8").Value = "B738"
$ws.Range("I8").Value = "(EI-EMH)"
$ws.Range("J8").Value = "8:49 PM"
$ws.Range("L8").Value = "0 hours, -21 minutes"

# K8 and M8 stay empty in the source table (same as K2:K7/M2:M7), but the
# sheet's used range still materializes them as blank cells. A plain
# Value="" assignment gets optimized away, so copy the (empty, default
# styled) formatting from A1 onto K8/M8 to force them to exist as real,
# default-style blank cells instead of being dropped on save.
$ws.Range("A1").Copy()
$ws.Range("K8").PasteSpecial(-4122)
$ws.Range("M8").PasteSpecial(-4122)
